# Update cryptocurrency price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.455.22'
$ws.Range('E2').Value = '  +1.84%  '

$ws.Range('D3').Value = '1.826.00'
$ws.Range('E3').Value = '  +1.66%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.27'

$ws.Range('E6').Value = '  -0.01%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5064'
$ws.Range('E7').Value = '  -4.69%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3905'
$ws.Range('E8').Value = '  +0.75%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07653'
$ws.Range('E9').Value = '  +2.70%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.84'
$ws.Range('E10').Value = '  +0.99%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.107'
$ws.Range('E11').Value = '  +1.83%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.07'
$ws.Range('E12').Value = '  +3.44%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.286'
$ws.Range('E13').Value = '  +1.79%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.585'
$ws.Range('E14').Value = '  +1.92%  '

$ws.Range('E15').Value = '  -0.02%  '

$ws.Range('D16').Value = '1.823.60'
$ws.Range('E16').Value = '  +1.75%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.09'
$ws.Range('E17').Value = '  +5.26%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001083'
$ws.Range('E18').Value = '  +2.16%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06673'
$ws.Range('E19').Value = '  +1.74%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.66'
$ws.Range('E20').Value = '  +2.33%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  +0.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.140'
$ws.Range('E22').Value = '  +3.03%  '

$ws.Range('D23').Value = '28.494.78'
$ws.Range('E23').Value = '  +1.91%  '

$ws.Range('E24').Value = '  +0.15%  '

$ws.Range('E25').Value = '  +7.77%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.67'
$ws.Range('E26').Value = '  -0.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.59'
$ws.Range('E27').Value = '  +2.19%  '

$ws.Range('D28').Value = '2.034.10'
$ws.Range('E28').Value = '  +1.75%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.394'
$ws.Range('E29').Value = '  +4.12%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.20'
$ws.Range('E30').Value = '  +2.65%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.126'
$ws.Range('E31').Value = '  +2.42%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1082'
$ws.Range('E32').Value = '  -0.62%  '

$ws.Range('E33').Value = '  +3.12%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.660'
$ws.Range('E34').Value = '  -0.18%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07037'
$ws.Range('E35').Value = '  +0.49%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2227'
$ws.Range('E36').Value = '  +0.98%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.952'
$ws.Range('E37').Value = '  +6.78%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02319'
$ws.Range('E38').Value = '  +1.92%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.132'
$ws.Range('E39').Value = '  +1.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6250'
$ws.Range('E40').Value = '  +2.17%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.21'
$ws.Range('E41').Value = '  -0.14%  '

$ws.Range('E42').Value = '  -0.76%  '

$ws.Range('E43').Value = '  -0.02%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.398'
$ws.Range('E44').Value = '  -1.47%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.43'
$ws.Range('E45').Value = '  +1.12%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5895'
$ws.Range('E46').Value = '  +3.28%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.718'
$ws.Range('E47').Value = '  +1.07%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.48'
$ws.Range('E48').Value = '  -0.26%  '

$ws.Range('E49').Value = '  +3.33%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.191'
$ws.Range('E50').Value = '  +1.10%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06914'
$ws.Range('E51').Value = '  +1.56%  '
